# Update the "想去人数" (interested-count) figures on the sheets that list
# individual events ("展览" and "全部类型") to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 1649
    $ws.Range("F6").Value = 631
    $ws.Range("F8").Value = 5747
}
